# Scheduled-runner update: refresh market-price-derived columns
# (currentAveragePrice / *NQ / *HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 9040
$ws.Range("I6").Value = 9040
$ws.Range("K6").Value = 27120
$ws.Range("M6").Value = -27008
# Row 16
$ws.Range("H16").Value = 4444
$ws.Range("I16").Value = 4444
$ws.Range("K16").Value = 4444
$ws.Range("M16").Value = -4214
# Row 86
$ws.Range("H86").Value = 4716.154
$ws.Range("I86").Value = 4059.1875
$ws.Range("K86").Value = 4059.1875
$ws.Range("M86").Value = -2936.1875
# Row 89
$ws.Range("H89").Value = 4716.154
$ws.Range("I89").Value = 4059.1875
$ws.Range("K89").Value = 20295.9375
$ws.Range("M89").Value = -14679.9375
# Row 113
$ws.Range("H113").Value = 3278.4285
$ws.Range("I113").Value = 2789.8
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 2789.8
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 464.1999999999998
$ws.Range("N113").Value = -11008
# Row 129
$ws.Range("H129").Value = 1107.8
$ws.Range("I129").Value = 887.25
$ws.Range("K129").Value = 2661.75
$ws.Range("M129").Value = 2338.25
# Row 131
$ws.Range("H131").Value = 2512.3845
$ws.Range("I131").Value = 923.7273
$ws.Range("J131").Value = 11250
$ws.Range("K131").Value = 2771.1819
$ws.Range("L131").Value = 33750
$ws.Range("M131").Value = 2268.8181
$ws.Range("N131").Value = -43830

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 749.6667
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 46
$ws.Range("H46").Value = 9045.666999999999
$ws.Range("J46").Value = 9045.666999999999
$ws.Range("L46").Value = 9045.666999999999
$ws.Range("N46").Value = -9683.666999999999
# Row 93
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992
# Row 94
$ws.Range("H94").Value = 30582.5
$ws.Range("J94").Value = 30582.5
$ws.Range("L94").Value = 30582.5
$ws.Range("N94").Value = -32384.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 749.6667
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 5
$ws.Range("H5").Value = 719
$ws.Range("I5").Value = 773.75
$ws.Range("K5").Value = 773.75
$ws.Range("M5").Value = -660.75
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
# Row 86
$ws.Range("H86").Value = 1999.4
$ws.Range("J86").Value = 1998
$ws.Range("L86").Value = 1998
$ws.Range("N86").Value = -4244
# Row 89
$ws.Range("H89").Value = 1999.4
$ws.Range("J89").Value = 1998
$ws.Range("L89").Value = 9990
$ws.Range("N89").Value = -21222

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 129.16667
$ws.Range("I7").Value = 75.5
$ws.Range("K7").Value = 75.5
$ws.Range("M7").Value = 37.5
# Row 28
$ws.Range("H28").Value = 10000
$ws.Range("J28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("N28").Value = -10490
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 100
$ws.Range("H100").Value = 99900
$ws.Range("J100").Value = 99900
$ws.Range("L100").Value = 99900
$ws.Range("N100").Value = -102064
# Row 130
$ws.Range("H130").Value = 63852
$ws.Range("J130").Value = 63852
$ws.Range("L130").Value = 63852
$ws.Range("N130").Value = -73892
# Row 132
$ws.Range("H132").Value = 1389.7646
$ws.Range("I132").Value = 1253.2903
$ws.Range("K132").Value = 3759.8709
$ws.Range("M132").Value = -1229.8709

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 34115280
$ws.Range("I11").Value = 37526730
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 112580190
$ws.Range("L11").Value = 2400
$ws.Range("M11").Value = -112580050
$ws.Range("N11").Value = -2680
# Row 80
$ws.Range("H80").Value = 2176.3333
$ws.Range("I80").Value = 2177
$ws.Range("K80").Value = 6531
$ws.Range("M80").Value = -5595
# Row 83
$ws.Range("H83").Value = 2176.3333
$ws.Range("I83").Value = 2177
$ws.Range("K83").Value = 19593
$ws.Range("M83").Value = -14913
# Row 92
$ws.Range("H92").Value = 251.44827
$ws.Range("I92").Value = 177.18182
$ws.Range("K92").Value = 531.5454599999999
$ws.Range("M92").Value = 716.4545400000001
# Row 122
$ws.Range("J122").Value = 99
$ws.Range("L122").Value = 891
$ws.Range("N122").Value = -5791
# Row 131
$ws.Range("H131").Value = 2935.7144
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 2935.7144
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 8807.143199999999
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -18887.1432

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 719999.5
$ws.Range("I7").Value = 719999.5
$ws.Range("K7").Value = 719999.5
$ws.Range("M7").Value = -719887.5
# Row 8
$ws.Range("H8").Value = 719999.5
$ws.Range("I8").Value = 719999.5
$ws.Range("K8").Value = 719999.5
$ws.Range("M8").Value = -719860.5
# Row 43
$ws.Range("H43").Value = 12613.6
$ws.Range("I43").Value = 12613.6
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 12613.6
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -12462.6
$ws.Range("N43").ClearContents()
# Row 102
$ws.Range("H102").Value = 1052.8334
$ws.Range("I102").Value = 1032.7142
$ws.Range("K102").Value = 1032.7142
$ws.Range("M102").Value = 589.2858000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3682.9167
$ws.Range("I40").Value = 3099.4285
$ws.Range("K40").Value = 3099.4285
$ws.Range("M40").Value = -2963.4285
# Row 122
$ws.Range("H122").Value = 7826.174
$ws.Range("I122").Value = 8469.385
$ws.Range("K122").Value = 25408.155
$ws.Range("M122").Value = -22958.155
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 49
$ws.Range("H49").Value = 44999
$ws.Range("J49").Value = 44999
$ws.Range("L49").Value = 44999
$ws.Range("N49").Value = -45459
# Row 51
$ws.Range("H51").Value = 41250
$ws.Range("I51").Value = 27500
$ws.Range("K51").Value = 27500
$ws.Range("M51").Value = -26990
# Row 75
$ws.Range("H75").Value = 90000
$ws.Range("I75").Value = 90000
$ws.Range("K75").Value = 90000
$ws.Range("M75").Value = -89064
# Row 78
$ws.Range("H78").Value = 90000
$ws.Range("I78").Value = 90000
$ws.Range("K78").Value = 270000
$ws.Range("M78").Value = -265320
# Row 81
$ws.Range("H81").Value = 2001047.8
$ws.Range("J81").Value = 5000400.5
$ws.Range("L81").Value = 10000801
$ws.Range("N81").Value = -10002923
# Row 84
$ws.Range("H84").Value = 2001047.8
$ws.Range("J84").Value = 5000400.5
$ws.Range("L84").Value = 50004005
$ws.Range("N84").Value = -50014613
# Row 100
$ws.Range("H100").Value = 11113882
$ws.Range("I100").Value = 12501867
$ws.Range("K100").Value = 25003734
$ws.Range("M100").Value = -25003193
# Row 126
$ws.Range("H126").Value = 4485.1
$ws.Range("I126").Value = 3013.7334
$ws.Range("K126").Value = 9041.200199999999
$ws.Range("M126").Value = -6571.200199999999
# Row 136
$ws.Range("H136").Value = 3348.75
$ws.Range("I136").Value = 3298.3333
$ws.Range("K136").Value = 9894.999899999999
$ws.Range("M136").Value = -7344.999899999999
